{"js": "const body = context.document.body;\n\n// The paragraph currently reads \"Version 2.\" split across four runs:\n//   \"Versi\" | \"on\" | \" 2\" | \".\"\n// (the first two runs are wrapped by a spell-check proofErr range, and a\n// _GoBack bookmark sits between the \" 2\" run and the final \".\" run).\n// The target text is \"Version 1.\" laid out as exactly two runs:\n//   \"Version\" | \" 1.\"\n//\n// Word.Range objects returned by search() point at a fixed position in the\n// document; once an earlier edit shifts content around, a range captured\n// before that edit can no longer be trusted. So we re-run search()\n// immediately before each mutation instead of reusing older range objects.\n\n// 1) Merge \"Versi\" + \"on\" into a single \"Version\" run: delete the \"on\" run,\n//    then rewrite \"Versi\" in place as \"Version\".\nlet results = body.search(\"on\", { matchCase: true });\nawait context.sync();\nresults.items[0].delete();\nawait context.sync();\n\nresults = body.search(\"Versi\", { matchCase: true });\nawait context.sync();\nresults.items[0].insertText(\"Version\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Merge \" 2\" + \".\" into a single \" 1.\" run: delete the trailing \".\" run,\n//    then rewrite \" 2\" in place as \" 1.\".\nresults = body.search(\".\", { matchCase: true });\nawait context.sync();\nresults.items[0].delete();\nawait context.sync();\n\nresults = body.search(\" 2\", { matchCase: true });\nawait context.sync();\nresults.items[0].insertText(\" 1.\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The paragraph currently reads \"Version 2.\" split across four runs:\n#   \"Versi\" | \"on\" | \" 2\" | \".\"\n# (the first two runs sit inside a spell-check proofErr range, and a\n# _GoBack bookmark sits between the \" 2\" run and the final \".\" run).\n# The target text is \"Version 1.\" laid out as exactly two runs:\n#   \"Version\" | \" 1.\"\n#\n# Each Find.Execute() narrows its Range to the match; that Range becomes\n# stale once a different part of the document is edited, so we grab a\n# fresh $d.Content range and re-run Find immediately before every edit\n# instead of reusing an older Range object.\n\n# 1) Merge \"Versi\" + \"on\" into a single \"Version\" run: delete the \"on\" run,\n#    then rewrite \"Versi\" in place as \"Version\".\n$r = $d.Content\n$r.Find.Text = \"on\"\n$r.Find.Execute() | Out-Null\n$r.Delete()\n\n$r = $d.Content\n$r.Find.Text = \"Versi\"\n$r.Find.Execute() | Out-Null\n$r.Text = \"Version\"\n\n# 2) Merge \" 2\" + \".\" into a single \" 1.\" run: delete the trailing \".\" run,\n#    then rewrite \" 2\" in place as \" 1.\".\n$r = $d.Content\n$r.Find.Text = \".\"\n$r.Find.Execute() | Out-Null\n$r.Delete()\n\n$r = $d.Content\n$r.Find.Text = \" 2\"\n$r.Find.Execute() | Out-Null\n$r.Text = \" 1.\"\n"}
